$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collapse Rare Categories")

# --- Copy the Action/Time/Content block formatting (F2:H7) into the new
#     "Data Polish" block (J2:L7) so the new columns inherit the same
#     bold-header / wrap-text formats already used for the other two
#     method blocks on this sheet. ---
$ws.Range("F2:H7").Copy()
$ws.Range("J2").PasteSpecial(-4122)

# Header row
$ws.Range("J2").Value = "Action"
$ws.Range("K2").Value = "Time"
$ws.Range("L2").Value = "Content"

# Row 3 - Load Data
$ws.Range("J3").Value = "Load Data"
$ws.Range("K3").Value = "1 min"
$ws.Range("L3").Value = "Upload the dataset to the 'Import' page."

# Row 4 - Inspect Data
$ws.Range("J4").Value = "Inspect Data"
$ws.Range("K4").Value = "3 min"
$ws.Range("L4").Value = "Examine category frequencies on the 'Data Profiling' page."

# Row 5 - Preprocess
$ws.Range("J5").Value = "Preprocess"
$ws.Range("K5").Value = "1 min"
$ws.Range("L5").Value = "Collapse rare categories on the 'Data Cleaning' page."

# Row 6 - Verify Changes
$ws.Range("J6").Value = "Verify Changes"
$ws.Range("K6").Value = "1 min"
$ws.Range("L6").Value = "Check collapsed categories on the 'Data Profiling' page."

# Row 7 - Overall / count, matches the "Overall" style used in B8:C8 /
# F8:G8 (bold, no wrap) rather than the wrapped body style above.
$ws.Range("B8:C8").Copy()
$ws.Range("J7").PasteSpecial(-4122)
$ws.Range("J7").Value = "Overall"
$ws.Range("K7").Value = 6
$ws.Range("L7").Font.Size = 12
$ws.Range("L7").WrapText = $true

# Row 1 trailing blank cell (matches D1/H1 pattern)
$ws.Range("H1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("L1").ClearContents()

# Column widths for the two new columns
$ws.Columns.Item(10).ColumnWidth = 14.1640625
$ws.Columns.Item(12).ColumnWidth = 24.83203125

# Make "Collapse Rare Categories" the active sheet/tab (this also clears
# tabSelected from the previously active "Remove Stopwords" sheet and
# updates the workbook's activeTab).
$ws.Activate()
$ws.Range("L1").Select()
$ws.Range("L1:L1048576").Select()

$wb.Windows.Item(1).ScrollWorkbookTabs(1, 2)
